$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new helper / query columns
# (order matters for the shared-string table: QUERY, NOMBRE, VOLUMEN, SABOR)
$ws.Range("J1").Value = "QUERY"
$ws.Range("F1").Value = "NOMBRE"
$ws.Range("G1").Value = "VOLUMEN"
$ws.Range("H1").Value = "SABOR"

# Row 2 gets the "source" formulas (not shared yet)
$ws.Range("F2").Formula = '=_xlfn.CONCAT(" ''",LEFT(B2,FIND(" -",B2,1)-1),"'', ")'
$ws.Range("G2").Formula = '=_xlfn.CONCAT("''",SUBSTITUTE(MID(SUBSTITUTE(" - " &B2&REPT(" ",6)," - ",REPT(",",255)),2*255,255),",",""),"'', ")'
$ws.Range("H2").Formula = '=_xlfn.CONCAT("''",TRIM(RIGHT(SUBSTITUTE(B2," ",REPT(" ",100)),100)),"'', ")'
$ws.Range("J2").Formula = '=_xlfn.CONCAT("INSERT INTO `tbproductos`(producto, nombre, envase, volumen, sabor, precio) VALUES (",A2,", ",F2,"''",C2,"'', ", G2,H2,D2, ");")'

# Rows 3-36: fill down as shared formulas (relative refs auto-adjust)
$ws.Range("F3:F36").Formula = '=_xlfn.CONCAT(" ''",LEFT(B3,FIND(" -",B3,1)-1),"'', ")'
$ws.Range("G3:G36").Formula = '=_xlfn.CONCAT("''",SUBSTITUTE(MID(SUBSTITUTE(" - " &B3&REPT(" ",6)," - ",REPT(",",255)),2*255,255),",",""),"'', ")'
$ws.Range("H3:H36").Formula = '=_xlfn.CONCAT("''",TRIM(RIGHT(SUBSTITUTE(B3," ",REPT(" ",100)),100)),"'', ")'
$ws.Range("J3:J36").Formula = '=_xlfn.CONCAT("INSERT INTO `tbproductos`(producto, nombre, envase, volumen, sabor, precio) VALUES (",A3,", ",F3,"''",C3,"'', ", G3,H3,D3, ");")'

# View / selection tweaks to mirror the author's final cursor position
$ws.Range("L23").Select()
